$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44294
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15000
$ws.Range("Q2").Value = '$/caja 16 kilos'
$ws.Range("R2").Value = 'Provincia de Los Andes'
$ws.Range("S2").Value = 938
$ws.Range("T2").Value = 16

# Row 3
$ws.Range("D3").Value = 44274
$ws.Range("L3").Value = 'Especial'
$ws.Range("M3").Value = 35
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 16000
$ws.Range("Q3").Value = '$/caja 16 kilos'
$ws.Range("R3").Value = 'Provincia de Los Andes'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 16

# Row 4
$ws.Range("D4").Value = 44274
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("Q4").Value = '$/caja 16 kilos'
$ws.Range("R4").Value = 'Provincia de Los Andes'
$ws.Range("S4").Value = 875
$ws.Range("T4").Value = 16

# Row 5
$ws.Range("D5").Value = 44258
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Provincia de Los Andes'
$ws.Range("S5").Value = 889
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44319
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("Q6").Value = '$/caja 16 kilos'
$ws.Range("R6").Value = 'Provincia de Los Andes'
$ws.Range("S6").Value = 938
$ws.Range("T6").Value = 16

# Row 7
$ws.Range("D7").Value = 44253
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 16000
$ws.Range("Q7").Value = '$/caja 16 kilos'
$ws.Range("R7").Value = 'Provincia de Los Andes'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 16

# Row 8
$ws.Range("D8").Value = 44418
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 30000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 30000
$ws.Range("Q8").Value = '$/caja 16 kilos'
$ws.Range("R8").Value = 'Provincia de Los Andes'
$ws.Range("S8").Value = 1875
$ws.Range("T8").Value = 16

# Row 9
$ws.Range("D9").Value = 44427
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 65
$ws.Range("N9").Value = 40000
$ws.Range("O9").Value = 40000
$ws.Range("P9").Value = 40000
$ws.Range("Q9").Value = '$/caja 16 kilos'
$ws.Range("R9").Value = 'Provincia de Los Andes'
$ws.Range("S9").Value = 2500
$ws.Range("T9").Value = 16

# Row 10
$ws.Range("D10").Value = 44280
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 20000
$ws.Range("P10").Value = 20000
$ws.Range("Q10").Value = '$/caja 16 kilos'
$ws.Range("R10").Value = 'Provincia de Los Andes'
$ws.Range("S10").Value = 1250
$ws.Range("T10").Value = 16

# Row 11
$ws.Range("D11").Value = 44221
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 25000
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 1389
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44267
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15400
$ws.Range("Q12").Value = '$/caja 16 kilos'
$ws.Range("R12").Value = 'Provincia de Los Andes'
$ws.Range("S12").Value = 962
$ws.Range("T12").Value = 16

# Row 13
$ws.Range("D13").Value = 44266
$ws.Range("L13").Value = 'Especial'
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 19000
$ws.Range("P13").Value = 18400
$ws.Range("Q13").Value = '$/caja 16 kilos'
$ws.Range("R13").Value = 'Provincia de Los Andes'
$ws.Range("S13").Value = 1150
$ws.Range("T13").Value = 16

# Row 14
$ws.Range("D14").Value = 44266
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 120
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("Q14").Value = '$/caja 16 kilos'
$ws.Range("R14").Value = 'Provincia de Los Andes'
$ws.Range("S14").Value = 938
$ws.Range("T14").Value = 16

# Row 15
$ws.Range("D15").Value = 44266
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("Q15").Value = '$/caja 16 kilos'
$ws.Range("R15").Value = 'Provincia de Los Andes'
$ws.Range("S15").Value = 625
$ws.Range("T15").Value = 16

# Row 16
$ws.Range("D16").Value = 44273
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 55
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("Q16").Value = '$/caja 16 kilos'
$ws.Range("R16").Value = 'Provincia de Los Andes'
$ws.Range("S16").Value = 875
$ws.Range("T16").Value = 16

# Row 17
$ws.Range("D17").Value = 44435
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 45000
$ws.Range("O17").Value = 45000
$ws.Range("P17").Value = 45000
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("R17").Value = 'Provincia de Los Andes'
$ws.Range("S17").Value = 2500
$ws.Range("T17").Value = 18

# Row 18
$ws.Range("D18").Value = 44277
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 90
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 16000
$ws.Range("Q18").Value = '$/caja 16 kilos'
$ws.Range("R18").Value = 'Provincia de Los Andes'
$ws.Range("S18").Value = 1000
$ws.Range("T18").Value = 16

# Row 19
$ws.Range("D19").Value = 44235
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 55
$ws.Range("N19").Value = 25000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 25000
$ws.Range("Q19").Value = '$/caja 16 kilos'
$ws.Range("R19").Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 1562
$ws.Range("T19").Value = 16

# Row 20
$ws.Range("D20").Value = 44244
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 14000
$ws.Range("Q20").Value = '$/caja 16 kilos'
$ws.Range("R20").Value = 'Provincia de Los Andes'
$ws.Range("S20").Value = 875
$ws.Range("T20").Value = 16

# Row 21
$ws.Range("D21").Value = 44320
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 14000
$ws.Range("O21").Value = 14000
$ws.Range("P21").Value = 14000
$ws.Range("Q21").Value = '$/caja 16 kilos'
$ws.Range("R21").Value = 'Provincia de Los Andes'
$ws.Range("S21").Value = 875
$ws.Range("T21").Value = 16

# Row 22
$ws.Range("D22").Value = 44326
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 15000
$ws.Range("Q22").Value = '$/caja 16 kilos'
$ws.Range("R22").Value = 'Provincia de Los Andes'
$ws.Range("S22").Value = 938
$ws.Range("T22").Value = 16

# Row 23
$ws.Range("D23").Value = 44321
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 15000
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 15000
$ws.Range("Q23").Value = '$/caja 16 kilos'
$ws.Range("R23").Value = 'Provincia de Los Andes'
$ws.Range("S23").Value = 938
$ws.Range("T23").Value = 16

# Row 24
$ws.Range("D24").Value = 44306
$ws.Range("L24").Value = 'Especial'
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 22000
$ws.Range("O24").Value = 22000
$ws.Range("P24").Value = 22000
$ws.Range("Q24").Value = '$/caja 16 kilos'
$ws.Range("R24").Value = 'Provincia de Los Andes'
$ws.Range("S24").Value = 1375
$ws.Range("T24").Value = 16

# Row 25
$ws.Range("D25").Value = 44265
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 16000
$ws.Range("O25").Value = 16000
$ws.Range("P25").Value = 16000
$ws.Range("Q25").Value = '$/caja 16 kilos'
$ws.Range("R25").Value = 'Provincia de Los Andes'
$ws.Range("S25").Value = 1000
$ws.Range("T25").Value = 16

# Row 26
$ws.Range("D26").Value = 44237
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 20
$ws.Range("N26").Value = 25000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 25000
$ws.Range("Q26").Value = '$/caja 16 kilos'
$ws.Range("R26").Value = 'Provincia de Limarí'
$ws.Range("S26").Value = 1562
$ws.Range("T26").Value = 16

# Row 27
$ws.Range("D27").Value = 44323
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 30
$ws.Range("N27").Value = 14000
$ws.Range("O27").Value = 14000
$ws.Range("P27").Value = 14000
$ws.Range("Q27").Value = '$/caja 16 kilos'
$ws.Range("R27").Value = 'Provincia de Los Andes'
$ws.Range("S27").Value = 875
$ws.Range("T27").Value = 16

# Row 28
$ws.Range("D28").Value = 44249
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 110
$ws.Range("N28").Value = 16000
$ws.Range("O28").Value = 16000
$ws.Range("P28").Value = 16000
$ws.Range("Q28").Value = '$/caja 16 kilos'
$ws.Range("R28").Value = 'Provincia de Los Andes'
$ws.Range("S28").Value = 1000
$ws.Range("T28").Value = 16

# Row 29
$ws.Range("D29").Value = 44270
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 55
$ws.Range("N29").Value = 10000
$ws.Range("O29").Value = 10000
$ws.Range("P29").Value = 10000
$ws.Range("Q29").Value = '$/caja 16 kilos'
$ws.Range("R29").Value = 'Provincia de Los Andes'
$ws.Range("S29").Value = 625
$ws.Range("T29").Value = 16

# Row 30
$ws.Range("D30").Value = 44301
$ws.Range("L30").Value = 'Especial'
$ws.Range("M30").Value = 55
$ws.Range("N30").Value = 22000
$ws.Range("O30").Value = 22000
$ws.Range("P30").Value = 22000
$ws.Range("Q30").Value = '$/caja 16 kilos'
$ws.Range("R30").Value = 'Provincia de Los Andes'
$ws.Range("S30").Value = 1375
$ws.Range("T30").Value = 16

# Row 31
$ws.Range("D31").Value = 44301
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 85
$ws.Range("N31").Value = 19000
$ws.Range("O31").Value = 19000
$ws.Range("P31").Value = 19000
$ws.Range("Q31").Value = '$/caja 16 kilos'
$ws.Range("R31").Value = 'Provincia de Los Andes'
$ws.Range("S31").Value = 1188
$ws.Range("T31").Value = 16

# Row 32
$ws.Range("D32").Value = 44284
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 25
$ws.Range("N32").Value = 18000
$ws.Range("O32").Value = 18000
$ws.Range("P32").Value = 18000
$ws.Range("Q32").Value = '$/caja 16 kilos'
$ws.Range("R32").Value = 'Provincia de Los Andes'
$ws.Range("S32").Value = 1125
$ws.Range("T32").Value = 16

# Row 33
$ws.Range("D33").Value = 44425
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 35
$ws.Range("N33").Value = 40000
$ws.Range("O33").Value = 40000
$ws.Range("P33").Value = 40000
$ws.Range("Q33").Value = '$/caja 16 kilos'
$ws.Range("R33").Value = 'Provincia de Los Andes'
$ws.Range("S33").Value = 2500
$ws.Range("T33").Value = 16

# Row 34
$ws.Range("D34").Value = 44315
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 105
$ws.Range("N34").Value = 15000
$ws.Range("O34").Value = 16000
$ws.Range("P34").Value = 15619
$ws.Range("Q34").Value = '$/caja 16 kilos'
$ws.Range("R34").Value = 'Provincia de Los Andes'
$ws.Range("S34").Value = 976
$ws.Range("T34").Value = 16

# Row 35
$ws.Range("D35").Value = 44327
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 35
$ws.Range("N35").Value = 15000
$ws.Range("O35").Value = 15000
$ws.Range("P35").Value = 15000
$ws.Range("Q35").Value = '$/caja 16 kilos'
$ws.Range("R35").Value = 'Provincia de Los Andes'
$ws.Range("S35").Value = 938
$ws.Range("T35").Value = 16

# Row 36
$ws.Range("D36").Value = 44238
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 65
$ws.Range("N36").Value = 14000
$ws.Range("O36").Value = 14000
$ws.Range("P36").Value = 14000
$ws.Range("Q36").Value = '$/caja 16 kilos'
$ws.Range("R36").Value = 'Provincia de Los Andes'
$ws.Range("S36").Value = 875
$ws.Range("T36").Value = 16

# Row 37
$ws.Range("D37").Value = 44236
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 45
$ws.Range("N37").Value = 25000
$ws.Range("O37").Value = 25000
$ws.Range("P37").Value = 25000
$ws.Range("Q37").Value = '$/caja 16 kilos'
$ws.Range("R37").Value = 'Provincia de Limarí'
$ws.Range("S37").Value = 1562
$ws.Range("T37").Value = 16

# Row 38
$ws.Range("D38").Value = 44389
$ws.Range("L38").Value = 'Especial'
$ws.Range("M38").Value = 30
$ws.Range("N38").Value = 30000
$ws.Range("O38").Value = 30000
$ws.Range("P38").Value = 30000
$ws.Range("Q38").Value = '$/caja 16 kilos'
$ws.Range("R38").Value = 'Provincia de Los Andes'
$ws.Range("S38").Value = 1875
$ws.Range("T38").Value = 16

# Row 39
$ws.Range("D39").Value = 44251
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 55
$ws.Range("N39").Value = 16000
$ws.Range("O39").Value = 16000
$ws.Range("P39").Value = 16000
$ws.Range("Q39").Value = '$/caja 16 kilos'
$ws.Range("R39").Value = 'Provincia de Los Andes'
$ws.Range("S39").Value = 1000
$ws.Range("T39").Value = 16

# Row 40
$ws.Range("D40").Value = 44330
$ws.Range("L40").Value = 'Especial'
$ws.Range("M40").Value = 35
$ws.Range("N40").Value = 30000
$ws.Range("O40").Value = 30000
$ws.Range("P40").Value = 30000
$ws.Range("Q40").Value = '$/caja 20 kilos'
$ws.Range("R40").Value = 'Provincia de Los Andes'
$ws.Range("S40").Value = 1500
$ws.Range("T40").Value = 20

# Row 41
$ws.Range("D41").Value = 44432
$ws.Range("L41").Value = 'Especial'
$ws.Range("M41").Value = 30
$ws.Range("N41").Value = 45000
$ws.Range("O41").Value = 45000
$ws.Range("P41").Value = 45000
$ws.Range("Q41").Value = '$/caja 18 kilos'
$ws.Range("R41").Value = 'Provincia de Los Andes'
$ws.Range("S41").Value = 2500
$ws.Range("T41").Value = 18
